$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B and D (header row + 4 data rows) - columns B and D
# trade places while column C stays fixed.
for ($row = 1; $row -le 5; $row++) {
    $bVal = $ws.Cells.Item($row, 2).Value()
    $dVal = $ws.Cells.Item($row, 4).Value()
    $ws.Cells.Item($row, 2).Value = $dVal
    $ws.Cells.Item($row, 4).Value = $bVal
}
